# Katalon Studio Bio-Farma automation test data update
# Third Commit (Ruli) - katalon bio farma v.01

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PR")

# Column C (Kode Material): was a bare number, now a text automation code
$ws.Range("C2").Value = "AutomatedTest001"
$ws.Range("C3").Value = "AutomatedTest002"

# Column D (Tanggal PR): refreshed dates
$ws.Range("D2").Value = "2024-01-04"
$ws.Range("D3").Value = "2024-02-04"

# Column E (NO PR): was a bare number, now a text automation code
$ws.Range("E2").Value = "M-AutomatedTest1"
$ws.Range("E3").Value = "M-AutomatedTest2"

# Column G (Spesifikasi): Otomatisai -> Otomatisasi (typo fix); G2 stays "Manual"
$ws.Range("G2").Value = "Manual"
$ws.Range("G3").Value = "Otomatisasi"

# Column O: unchanged value, rewritten for stability
$ws.Range("O3").Value = "TAKEHARA KAGAKU"

# Column P (Country): INDIA -> VATIKAN; P3 stays "JAPAN"
$ws.Range("P2").Value = "VATIKAN"
$ws.Range("P3").Value = "JAPAN"

# Column R (No Kontrak JP): R2 stays "JP-SAMPLE00018"; R3 newly populated
$ws.Range("R2").Value = "JP-SAMPLE00018"
$ws.Range("R3").Value = "JP-SAMPLE00019"

# Column S (Status PR): unchanged values
$ws.Range("S2").Value = "Approve"
$ws.Range("S3").Value = "Approve"

# Update the active selection to match the saved workbook state
$ws.Range("S5").Select()
